$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 4.75
$ws.Range("L6").Value = 5.5
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("X6").Value = 7.5
$ws.Range("Z6").Value = 17
$ws.Range("AI6").Value = 21
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 34
$ws.Range("AS6").Value = 401
$ws.Range("AX6").Value = 6
$ws.Range("AY6").Value = 29

# Row 7
$ws.Range("G7").Value = 2.75
$ws.Range("I7").Value = 2.75
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2.35
$ws.Range("R7").Value = 1.57
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.8
$ws.Range("AA7").Value = 23
$ws.Range("AC7").Value = 7.5
$ws.Range("AE7").Value = 15
$ws.Range("AL7").Value = 23

# Row 8
$ws.Range("G8").Value = 2.25
$ws.Range("H8").Value = 2.63
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 3.25
$ws.Range("Q8").Value = 4.2
$ws.Range("R8").Value = 1.22
$ws.Range("S8").Value = 1.95
$ws.Range("T8").Value = 1.85
$ws.Range("X8").Value = 8
$ws.Range("Z8").Value = 21
$ws.Range("AI8").Value = 19
$ws.Range("AO8").Value = 17

# Row 9
$ws.Range("G9").Value = 1.67
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 6.25
$ws.Range("J9").Value = 2.38
$ws.Range("L9").Value = 7
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("U9").Value = 2.5
$ws.Range("V9").Value = 1.5
$ws.Range("W9").Value = 4.75
$ws.Range("X9").Value = 6
$ws.Range("Z9").Value = 12
$ws.Range("AA9").Value = 17
$ws.Range("AC9").Value = 6
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 26
$ws.Range("AF9").Value = 101
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 21
$ws.Range("AK9").Value = 67
$ws.Range("AM9").Value = 67
$ws.Range("AN9").Value = 3.4
$ws.Range("AO9").Value = 9
$ws.Range("AP9").Value = 29
$ws.Range("AS9").Value = 301
$ws.Range("AU9").Value = 11
$ws.Range("AV9").Value = 101
$ws.Range("AX9").Value = 7.5
$ws.Range("AY9").Value = 41
$ws.Range("AZ9").Value = 51
$ws.Range("BA9").Value = 151
$ws.Range("BB9").Value = 251

# Row 12
$ws.Range("G12").Value = 1.83
$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 4.2
$ws.Range("J12").Value = 2.4
$ws.Range("L12").Value = 4.33
$ws.Range("U12").Value = 1.67
$ws.Range("V12").Value = 2.1
$ws.Range("W12").Value = 8.5
$ws.Range("X12").Value = 9.5
$ws.Range("Z12").Value = 15
$ws.Range("AC12").Value = 13
$ws.Range("AD12").Value = 7
$ws.Range("AE12").Value = 13
$ws.Range("AF12").Value = 41
$ws.Range("AG12").Value = 151
$ws.Range("AH12").Value = 13
$ws.Range("AI12").Value = 21
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 41
$ws.Range("AL12").Value = 29
$ws.Range("AM12").Value = 34
$ws.Range("AN12").Value = 4
$ws.Range("AO12").Value = 9.5
$ws.Range("AQ12").Value = 29
$ws.Range("AR12").Value = 51
$ws.Range("AU12").Value = 7.5
$ws.Range("AX12").Value = 6
$ws.Range("AY12").Value = 21
$ws.Range("AZ12").Value = 26
$ws.Range("BA12").Value = 67
$ws.Range("BB12").Value = 81
$ws.Range("BC12").Value = 151
